# Update "想去人数" (want-to-go count) figures with freshly scraped numbers.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 8198
$wsExpo.Range("F5").Value = 5969
$wsExpo.Range("F6").Value = 507
$wsExpo.Range("F11").Value = 761
$wsExpo.Range("F12").Value = 76

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 93

# Sheet "全部类型" (All types, the combined listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 8198
$wsAll.Range("F5").Value = 5969
$wsAll.Range("F6").Value = 507
$wsAll.Range("F11").Value = 93
$wsAll.Range("F15").Value = 761
$wsAll.Range("F16").Value = 76
